# Fixing RAD Test Cases and Data
# Update the test execution timestamps in column B (rows 2-4) with new run dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "Fri Oct 21 17:41:34 EDT 2022"
$ws.Range("B3").Value = "Fri Oct 21 17:41:46 EDT 2022"
$ws.Range("B4").Value = "Fri Oct 21 17:41:56 EDT 2022"
